# Bitacora Laboratorio - actualizacion 11 de febrero de 2024
# Adds P6/P7 encuadre-montaje-reporte columns (and Puntaje/Calificacion)
# to the "Concentrado" sheet, plus a few new "falta" marks on "Faltas".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Faltas")
$ws2 = $wb.Worksheets.Item("Concentrado")

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1) Faltas sheet: four new attendance marks in column G
# ---------------------------------------------------------------------
$ws1.Range("G2").Value = 1
$ws1.Range("G5").Value = 1
$ws1.Range("G10").Value = 1
$ws1.Range("G14").Value = 1

# ---------------------------------------------------------------------
# 2) Concentrado sheet: new header labels.
#    Written in this precise order so the shared-string table grows in
#    the same sequence the workbook expects (P6_Encuadre, P6_Montaje,
#    P5_Reporte, P6_Reporte, Puntaje, Calificacion, P5_Marco_Teorico,
#    P7_Encuadre, P7_Montaje).
# ---------------------------------------------------------------------
$ws2.Range("J1").Value = "P6_Encuadre"
$ws2.Range("K1").Value = "P6_Montaje"
$ws2.Range("I1").Value = "P5_Reporte"
$ws2.Range("L1").Value = "P6_Reporte"
$ws2.Range("O1").Value = "Puntaje"
$ws2.Range("P1").Value = "Calificación"
$ws2.Range("G1").Value = "P5_Marco_Teórico"
$ws2.Range("M1").Value = "P7_Encuadre"
$ws2.Range("N1").Value = "P7_Montaje"
$ws2.Range("H1").Value = "P5_Montaje"

$ws2.Range("G1:P1").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 3) Concentrado sheet: per-student scores.
#    Columns: row, E, G, H, I, J, K  (F already holds its final value)
# ---------------------------------------------------------------------
$rows = @(
    @(2, 0, 0, 5, 0, 5, 0),
    @(3, 4, 4.5, 5, 4, 5, 5),
    @(4, 2.4, 3.2, 5, 0, 5, 5),
    @(5, 2.5, 4, 5, 3, 5, 0),
    @(6, 4, 5, 5, 5, 5, 5),
    @(7, 0, 0, 5, 2.5, 5, 5),
    @(8, 4, 5, 5, 5, 5, 5),
    @(9, 5, 5, 5, 5, 5, 5),
    @(10, 0, 0, 5, 2.5, 5, 0),
    @(11, 2.5, 5, 5, 0, 5, 5),
    @(12, 0, 0, 5, 0, 5, 5),
    @(13, 4, 5, 5, 5, 5, 5),
    @(14, 4, 5, 5, 5, 5, 0),
    @(15, 5, 5, 5, 5, 5, 5)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws2.Cells.Item($r, 5).Value  = $row[1]   # E - P4_Reporte
    $ws2.Cells.Item($r, 7).Value  = $row[2]   # G - P5_Marco_Teorico
    $ws2.Cells.Item($r, 8).Value  = $row[3]   # H - P5_Montaje
    $ws2.Cells.Item($r, 9).Value  = $row[4]   # I - P5_Reporte
    $ws2.Cells.Item($r, 10).Value = $row[5]   # J - P6_Encuadre
    $ws2.Cells.Item($r, 11).Value = $row[6]   # K - P6_Montaje

    $ws2.Cells.Item($r, 11).HorizontalAlignment = $xlCenter

    $ws2.Cells.Item($r, 15).Formula = "=SUM(E" + $r + ":M" + $r + ")"
    $ws2.Cells.Item($r, 15).HorizontalAlignment = $xlCenter
}

# ---------------------------------------------------------------------
# 4) Column widths (best-fit approximation for the new columns)
# ---------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 6.85546875 - (5/6)
$ws1.Columns.Item(6).ColumnWidth = 6.85546875 - (5/6)
$ws1.Columns.Item(7).ColumnWidth = 6.85546875 - (5/6)
$ws1.Columns.Item(8).ColumnWidth = 6.85546875 - (5/6)

$ws2.Columns.Item(7).ColumnWidth  = 17.140625 - (5/6)
$ws2.Columns.Item(8).ColumnWidth  = 11.140625 - (5/6)
$ws2.Columns.Item(9).ColumnWidth  = 11.140625 - (5/6)
$ws2.Columns.Item(10).ColumnWidth = 12.42578125 - (5/6)
$ws2.Columns.Item(11).ColumnWidth = 11.140625 - (5/6)
$ws2.Columns.Item(13).ColumnWidth = 12.42578125 - (5/6)
$ws2.Columns.Item(14).ColumnWidth = 11.140625 - (5/6)
$ws2.Columns.Item(16).ColumnWidth = 11.85546875 - (5/6)

# ---------------------------------------------------------------------
# 5) Selection bookmarks (each sheet keeps its own last-used cell, and
#    Concentrado stays the visible/active tab, matching the source file)
# ---------------------------------------------------------------------
$ws1.Range("H2").Select()
$ws2.Range("L9").Select()
$ws2.Select()
